$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- Simple text replacements (paragraph count unchanged) ---

Replace-Text "Generated: 2026-02-09 14:18" "Generated: 2026-02-09 14:44"

Replace-Text "This application provides an anomaly inspection pipeline with two UIs: a user camera interface and an admin upload interface." `
             "The system provides a user camera UI and an admin upload UI to run anomaly inspection."

Replace-Text "2.1 src/api.py ? Web API + UI" "2.1 src/api.py — Web API + UI"

# "Routes:" paragraph absorbs the bullet list text, the bullets themselves
# are deleted as whole paragraphs further below.
Replace-Text "Routes:" "Routes: GET /, GET /admin, POST /upload, POST /analyze."

Replace-Text "2.2 src/pipeline.py ? Batch Pipeline" "2.2 src/pipeline.py — Batch Pipeline"
Replace-Text "Purpose: Process dataset test images end-to-end." "Batch processing for dataset evaluation."

Replace-Text "2.3 src/data/mvtec.py ? Dataset Loader" "2.3 src/data/mvtec.py — Dataset Loader"
Replace-Text "Purpose: Iterate the MVTec AD folder structure and link ground-truth masks." `
             "Loads MVTec AD style data and ground-truth masks."

Replace-Text "2.4 src/models/mean_diff.py ? Baseline Model" "2.4 src/models/mean_diff.py — Baseline Model"
Replace-Text "Purpose: Simple anomaly detection using mean and std from good images." `
             "Mean/std baseline model for anomaly scoring."

Replace-Text "2.5 src/postproc/* ? Heatmap, Mask, BBoxes" "2.5 src/postproc/* — Heatmap, Mask, BBoxes"
Replace-Text "Purpose: Normalize heatmaps, threshold to mask, and derive bounding boxes." `
             "Normalize, threshold, and bounding box extraction."

Replace-Text "2.6 src/vlm/semantics.py ? VLM Semantics (Placeholder)" "2.6 src/vlm/semantics.py — VLM"
Replace-Text "Purpose: Defect label + evidence from VLM. Currently returns Unknown." `
             "LLaVA-1.6 (Mistral) integration for defect labels."

Replace-Text "2.7 src/risk/* ? Risk Lookup + Policy" "2.7 src/risk/* — Risk Lookup + Policy"
Replace-Text "Purpose: RPM lookup and risk-to-action mapping (placeholders until RPM provided)." `
             "RPM lookup and risk-to-action mapping."

Replace-Text "2.8 src/uncertainty/* ? Confidence + Human Review" "2.8 src/uncertainty/* — Confidence + Review"
Replace-Text "Purpose: Combine confidence signals and decide if human review is required." `
             "Confidence fusion and human review rules."

# Old "2.9 src/report/writer.py" heading is repurposed into the new "2.9
# configs/base.yaml" heading (its old Heading2 style/paragraph is kept);
# the old 2.10 heading paragraph is removed separately below.
Replace-Text "2.9 src/report/writer.py ? Output Writer" "2.9 configs/base.yaml — Configuration"

Replace-Text "Purpose: Central configuration for categories, labels, thresholds, RPM, and paths." `
             "Central configuration for categories and thresholds."

Replace-Text "User UI flow:" "User UI flow: start camera, capture, analyze, review result."

Replace-Text "1) Upload an image with category and description." `
             "Admin UI flow: upload image, analyze, review result."

Replace-Text "- JSON results with anomaly score, decision, bboxes, and placeholders for risk/labels." `
             "Heatmaps, masks, and JSON outputs saved under artifacts/ and outputs/."

Replace-Text "5. Current Limitations" "5. Running"

Replace-Text "- Baseline model is simple; replace with PatchCore/PaDiM for higher accuracy." `
             "Install dependencies and run: python -m uvicorn src.api:app --reload"

# --- Paragraph deletions (whole paragraphs removed) ---
# Delete from bottom to top so earlier paragraph indices stay valid.

$deletions = @(
    52, # python -m uvicorn src.api:app --reload
    51, # Start server:
    50, # python -m pip install -r requirements.txt
    49, # Install dependencies:
    48, # 6. Running the Application
    47, # - RPM table and label sets must be populated.
    46, # - VLM integration is a placeholder.
    42, # - artifacts/masks/*_mask.png
    41, # - artifacts/heatmaps/*_hm.png
    39, # 2) Server analyzes and returns results.
    37, # Admin UI flow:
    36, # 3) Server analyzes and returns results.
    35, # 2) Submit image to server.
    34, # 1) Start camera, capture frame.
    30, # 2.10 configs/base.yaml ? Configuration
    29, # Purpose: Persist JSON outputs.
    15, # Process: Train baseline model on good images, infer on test images, save heatmaps/masks, write JSON results.
    12, # - POST /analyze : Save image, run anomaly analysis, return result
    11, # - POST /upload : Save uploaded image
    10, # - GET /admin : Admin upload UI
    9,  # - GET / : User camera UI
    7   # Purpose: HTTP endpoints and HTML interfaces.
)

foreach ($idx in $deletions) {
    $d.Paragraphs($idx).Range.Delete()
}
